$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 15).Value = 1.27
$ws.Cells.Item(2, 17).Value = 1.27
$ws.Cells.Item(2, 18).Value = 1.16
$ws.Cells.Item(2, 19).Value = 1.27
$ws.Cells.Item(3, 6).Value = 1.66
$ws.Cells.Item(3, 7).Value = 1.84
$ws.Cells.Item(3, 8).Value = 5
$ws.Cells.Item(3, 9).Value = 7.2
$ws.Cells.Item(3, 10).Value = 3.35
$ws.Cells.Item(3, 11).Value = 3.85
$ws.Cells.Item(3, 12).Value = 1.48
$ws.Cells.Item(3, 13).Value = 1.08
$ws.Cells.Item(3, 14).Value = 2.92
$ws.Cells.Item(3, 15).Value = 1.42
$ws.Cells.Item(3, 16).Value = 1.66
$ws.Cells.Item(3, 17).Value = 2.24
$ws.Cells.Item(3, 18).Value = 1.24
$ws.Cells.Item(3, 19).Value = 4.3
$ws.Cells.Item(3, 20).Value = 2.04
$ws.Cells.Item(3, 21).Value = 1.78
$ws.Cells.Item(3, 22).Value = 1.16
$ws.Cells.Item(3, 23).Value = 1.92
$ws.Cells.Item(3, 24).Value = 13.5
$ws.Cells.Item(3, 25).Value = 18.5
$ws.Cells.Item(3, 26).Value = 55
$ws.Cells.Item(3, 28).Value = 8.199999999999999
$ws.Cells.Item(3, 29).Value = 9.800000000000001
$ws.Cells.Item(3, 30).Value = 27
$ws.Cells.Item(3, 32).Value = 12
$ws.Cells.Item(3, 33).Value = 12.5
$ws.Cells.Item(3, 34).Value = 29
$ws.Cells.Item(3, 36).Value = 24
$ws.Cells.Item(3, 37).Value = 24
$ws.Cells.Item(3, 40).Value = 19.5
$ws.Cells.Item(4, 8).Value = 7.8
$ws.Cells.Item(4, 14).Value = 5.3
$ws.Cells.Item(4, 15).Value = 1.2
$ws.Cells.Item(4, 19).Value = 2.38
$ws.Cells.Item(4, 20).Value = 1.78
$ws.Cells.Item(4, 26).Value = 95
$ws.Cells.Item(4, 30).Value = 980
$ws.Cells.Item(4, 34).Value = 28
$ws.Cells.Item(4, 36).Value = 15.5
$ws.Cells.Item(4, 37).Value = 16.5
$ws.Cells.Item(4, 38).Value = 34
$ws.Cells.Item(4, 40).Value = 5.7
$ws.Cells.Item(5, 6).Value = 3.45
$ws.Cells.Item(5, 7).Value = 4.4
$ws.Cells.Item(5, 8).Value = 1.95
$ws.Cells.Item(5, 9).Value = 2.2
$ws.Cells.Item(5, 11).Value = 4.7
$ws.Cells.Item(5, 12).Value = 1.31
$ws.Cells.Item(5, 13).Value = 1.05
$ws.Cells.Item(5, 14).Value = 3.45
$ws.Cells.Item(5, 15).Value = 1.28
$ws.Cells.Item(5, 16).Value = 1.98
$ws.Cells.Item(5, 17).Value = 1.82
$ws.Cells.Item(5, 18).Value = 1.38
$ws.Cells.Item(5, 19).Value = 2.78
$ws.Cells.Item(5, 20).Value = 1.7
$ws.Cells.Item(5, 21).Value = 2.12
$ws.Cells.Item(5, 22).Value = 1.84
$ws.Cells.Item(5, 23).Value = 1.29
$ws.Cells.Item(6, 6).Value = 5.7
$ws.Cells.Item(6, 7).Value = 6.8
$ws.Cells.Item(6, 8).Value = 1.51
$ws.Cells.Item(6, 10).Value = 4.6
$ws.Cells.Item(6, 11).Value = 5.5
$ws.Cells.Item(6, 12).Value = 1.22
$ws.Cells.Item(6, 14).Value = 5.8
$ws.Cells.Item(6, 15).Value = 1.17
$ws.Cells.Item(6, 16).Value = 2.64
$ws.Cells.Item(6, 17).Value = 1.5
$ws.Cells.Item(6, 18).Value = 1.66
$ws.Cells.Item(6, 19).Value = 2.2
$ws.Cells.Item(6, 20).Value = 1.64
$ws.Cells.Item(6, 21).Value = 2.28
$ws.Cells.Item(6, 23).Value = 1.17
$ws.Cells.Item(6, 24).Value = 36
$ws.Cells.Item(6, 25).Value = 15.5
$ws.Cells.Item(6, 26).Value = 14.5
$ws.Cells.Item(6, 27).Value = 18.5
$ws.Cells.Item(6, 28).Value = 36
$ws.Cells.Item(6, 29).Value = 14.5
$ws.Cells.Item(6, 30).Value = 13
$ws.Cells.Item(6, 31).Value = 18
$ws.Cells.Item(6, 33).Value = 29
$ws.Cells.Item(6, 34).Value = 23
$ws.Cells.Item(6, 35).Value = 28
$ws.Cells.Item(6, 37).Value = 85
$ws.Cells.Item(6, 39).Value = 90
$ws.Cells.Item(6, 41).Value = 6.8
$ws.Cells.Item(7, 6).Value = 3.9
$ws.Cells.Item(7, 7).Value = 6
$ws.Cells.Item(7, 8).Value = 1.81
$ws.Cells.Item(7, 9).Value = 2.06
$ws.Cells.Item(7, 10).Value = 3
$ws.Cells.Item(7, 11).Value = 4.5
$ws.Cells.Item(7, 12).Value = 1.4
$ws.Cells.Item(7, 13).Value = 1.07
$ws.Cells.Item(7, 14).Value = 3.1
$ws.Cells.Item(7, 15).Value = 1.36
$ws.Cells.Item(7, 16).Value = 1.71
$ws.Cells.Item(7, 17).Value = 2.06
$ws.Cells.Item(7, 18).Value = 1.27
$ws.Cells.Item(7, 19).Value = 3.45
$ws.Cells.Item(7, 20).Value = 1.9
$ws.Cells.Item(7, 21).Value = 1.83
$ws.Cells.Item(7, 22).Value = 1.94
$ws.Cells.Item(7, 23).Value = 1.22
$ws.Cells.Item(8, 6).Value = 2.9
$ws.Cells.Item(8, 7).Value = 3.45
$ws.Cells.Item(8, 8).Value = 2.24
$ws.Cells.Item(8, 9).Value = 2.5
$ws.Cells.Item(8, 10).Value = 3.35
$ws.Cells.Item(8, 11).Value = 4.9
$ws.Cells.Item(8, 12).Value = 1.27
$ws.Cells.Item(8, 13).Value = 1.04
$ws.Cells.Item(8, 14).Value = 4
$ws.Cells.Item(8, 15).Value = 1.23
$ws.Cells.Item(8, 16).Value = 2.22
$ws.Cells.Item(8, 17).Value = 1.65
$ws.Cells.Item(8, 18).Value = 1.49
$ws.Cells.Item(8, 19).Value = 2.44
$ws.Cells.Item(8, 20).Value = 1.58
$ws.Cells.Item(8, 21).Value = 2.36
$ws.Cells.Item(8, 22).Value = 1.67
$ws.Cells.Item(8, 23).Value = 1.42
$ws.Cells.Item(9, 6).Value = 1.5
$ws.Cells.Item(9, 7).Value = 1.63
$ws.Cells.Item(9, 8).Value = 5.9
$ws.Cells.Item(9, 9).Value = 9.6
$ws.Cells.Item(9, 10).Value = 3.65
$ws.Cells.Item(9, 11).Value = 5.5
$ws.Cells.Item(9, 12).Value = 1.32
$ws.Cells.Item(9, 13).Value = 1.05
$ws.Cells.Item(9, 14).Value = 3.35
$ws.Cells.Item(9, 15).Value = 1.29
$ws.Cells.Item(9, 16).Value = 1.95
$ws.Cells.Item(9, 17).Value = 1.85
$ws.Cells.Item(9, 18).Value = 1.36
$ws.Cells.Item(9, 19).Value = 2.88
$ws.Cells.Item(9, 20).Value = 1.94
$ws.Cells.Item(9, 21).Value = 1.84
$ws.Cells.Item(9, 23).Value = 2.58
$ws.Cells.Item(10, 6).Value = 1.7
$ws.Cells.Item(10, 7).Value = 1.99
$ws.Cells.Item(10, 8).Value = 3.45
$ws.Cells.Item(10, 9).Value = 4.7
$ws.Cells.Item(10, 10).Value = 3.7
$ws.Cells.Item(10, 11).Value = 5
$ws.Cells.Item(10, 15).Value = 1.17
$ws.Cells.Item(10, 16).Value = 2.32
$ws.Cells.Item(10, 17).Value = 1.5
$ws.Cells.Item(10, 18).Value = 1.6
$ws.Cells.Item(10, 19).Value = 2.1
$ws.Cells.Item(10, 20).Value = 1.53
$ws.Cells.Item(10, 21).Value = 2.42
$ws.Cells.Item(10, 22).Value = 1.27
$ws.Cells.Item(10, 23).Value = 2
$ws.Cells.Item(10, 24).Value = 32
$ws.Cells.Item(10, 25).Value = 28
$ws.Cells.Item(10, 26).Value = 44
$ws.Cells.Item(10, 28).Value = 17
$ws.Cells.Item(10, 29).Value = 13
$ws.Cells.Item(10, 30).Value = 21
$ws.Cells.Item(10, 31).Value = 50
$ws.Cells.Item(10, 32).Value = 17.5
$ws.Cells.Item(10, 33).Value = 13
$ws.Cells.Item(10, 34).Value = 21
$ws.Cells.Item(10, 35).Value = 55
$ws.Cells.Item(10, 36).Value = 26
$ws.Cells.Item(10, 37).Value = 21
$ws.Cells.Item(10, 38).Value = 32
$ws.Cells.Item(10, 40).Value = 9.4
$ws.Cells.Item(10, 41).Value = 36
$ws.Cells.Item(11, 6).Value = 1.62
$ws.Cells.Item(11, 7).Value = 1.77
$ws.Cells.Item(11, 8).Value = 5.5
$ws.Cells.Item(11, 9).Value = 8.800000000000001
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(11, 11).Value = 4.6
$ws.Cells.Item(11, 12).Value = 1.42
$ws.Cells.Item(11, 13).Value = 1.1
$ws.Cells.Item(11, 14).Value = 2.52
$ws.Cells.Item(11, 15).Value = 1.45
$ws.Cells.Item(11, 16).Value = 1.6
$ws.Cells.Item(11, 17).Value = 2.14
$ws.Cells.Item(11, 18).Value = 1.22
$ws.Cells.Item(11, 19).Value = 4.1
$ws.Cells.Item(11, 20).Value = 2.02
$ws.Cells.Item(11, 21).Value = 1.69
$ws.Cells.Item(11, 22).Value = 1.14
$ws.Cells.Item(11, 23).Value = 2.3
$ws.Cells.Item(12, 6).Value = 2.86
$ws.Cells.Item(12, 7).Value = 3.1
$ws.Cells.Item(12, 8).Value = 2.8
$ws.Cells.Item(12, 9).Value = 3.05
$ws.Cells.Item(12, 10).Value = 3.1
$ws.Cells.Item(12, 11).Value = 3.2
$ws.Cells.Item(12, 16).Value = 1.57
$ws.Cells.Item(12, 17).Value = 2.5
$ws.Cells.Item(12, 22).Value = 1.5
$ws.Cells.Item(12, 23).Value = 1.48
$ws.Cells.Item(12, 28).Value = 10
$ws.Cells.Item(12, 29).Value = 8.6
$ws.Cells.Item(12, 32).Value = 22
$ws.Cells.Item(12, 33).Value = 980
$ws.Cells.Item(13, 6).Value = 1.96
$ws.Cells.Item(13, 7).Value = 2.1
$ws.Cells.Item(13, 8).Value = 4.5
$ws.Cells.Item(13, 9).Value = 5.1
$ws.Cells.Item(13, 10).Value = 3.25
$ws.Cells.Item(13, 11).Value = 3.55
$ws.Cells.Item(13, 12).Value = 1.51
$ws.Cells.Item(13, 13).Value = 1.1
$ws.Cells.Item(13, 14).Value = 2.9
$ws.Cells.Item(13, 15).Value = 1.43
$ws.Cells.Item(13, 16).Value = 1.64
$ws.Cells.Item(13, 17).Value = 2.34
$ws.Cells.Item(13, 18).Value = 1.22
$ws.Cells.Item(13, 19).Value = 4.6
$ws.Cells.Item(13, 20).Value = 2.06
$ws.Cells.Item(13, 21).Value = 1.76
$ws.Cells.Item(13, 22).Value = 1.24
$ws.Cells.Item(13, 23).Value = 1.9
$ws.Cells.Item(13, 24).Value = 10.5
$ws.Cells.Item(13, 27).Value = 150
$ws.Cells.Item(13, 28).Value = 8.800000000000001
$ws.Cells.Item(13, 29).Value = 9.199999999999999
$ws.Cells.Item(13, 30).Value = 25
$ws.Cells.Item(13, 32).Value = 13
$ws.Cells.Item(13, 33).Value = 14
$ws.Cells.Item(13, 35).Value = 120
$ws.Cells.Item(13, 38).Value = 65
$ws.Cells.Item(13, 39).Value = 210
$ws.Cells.Item(13, 40).Value = 25
$ws.Cells.Item(13, 41).Value = 130
$ws.Cells.Item(14, 6).Value = 2.22
$ws.Cells.Item(14, 7).Value = 2.38
$ws.Cells.Item(14, 8).Value = 3.65
$ws.Cells.Item(14, 9).Value = 4.1
$ws.Cells.Item(14, 10).Value = 3.1
$ws.Cells.Item(14, 11).Value = 3.45
$ws.Cells.Item(14, 12).Value = 1.55
$ws.Cells.Item(14, 14).Value = 2.7
$ws.Cells.Item(14, 15).Value = 1.5
$ws.Cells.Item(14, 16).Value = 1.56
$ws.Cells.Item(14, 17).Value = 2.46
$ws.Cells.Item(14, 19).Value = 4.8
$ws.Cells.Item(14, 20).Value = 2.04
$ws.Cells.Item(14, 21).Value = 1.76
$ws.Cells.Item(14, 22).Value = 1.32
$ws.Cells.Item(14, 23).Value = 1.72
$ws.Cells.Item(14, 29).Value = 8
$ws.Cells.Item(14, 38).Value = 75
$ws.Cells.Item(15, 6).Value = 2.3
$ws.Cells.Item(15, 7).Value = 2.5
$ws.Cells.Item(15, 8).Value = 3.55
$ws.Cells.Item(15, 9).Value = 4.1
$ws.Cells.Item(15, 10).Value = 3.05
$ws.Cells.Item(15, 11).Value = 3.15
$ws.Cells.Item(15, 13).Value = 1.12
$ws.Cells.Item(15, 14).Value = 2.62
$ws.Cells.Item(15, 16).Value = 1.54
$ws.Cells.Item(15, 17).Value = 2.58
$ws.Cells.Item(15, 20).Value = 2.12
$ws.Cells.Item(15, 21).Value = 1.76
$ws.Cells.Item(15, 22).Value = 1.32
$ws.Cells.Item(15, 23).Value = 1.66
$ws.Cells.Item(15, 26).Value = 28
$ws.Cells.Item(15, 27).Value = 90
$ws.Cells.Item(15, 28).Value = 8.800000000000001
$ws.Cells.Item(15, 29).Value = 7.6
$ws.Cells.Item(15, 30).Value = 18
$ws.Cells.Item(15, 31).Value = 70
$ws.Cells.Item(15, 32).Value = 14
$ws.Cells.Item(15, 33).Value = 13
$ws.Cells.Item(15, 37).Value = 40
$ws.Cells.Item(15, 38).Value = 80
$ws.Cells.Item(15, 39).Value = 240
$ws.Cells.Item(15, 40).Value = 36
